$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-18 in column C ("Förändrad") have their date serial value
# bumped from 45179 (2023-09-10) to 45180 (2023-09-11).
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45180
}
